# Auto-generated edit script applying the diff changes to CryCompanywiseStockReport_1.xlsx
# Sets each changed cell to its new literal value (sheet contains static values, no formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("F52").Value = 34
$ws.Range("G52").Value = 1516.4
$ws.Range("B61").Value = 25393.2
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("B176").Value = 15054.45
$ws.Range("B181").Value = 53925
$ws.Range("B182").Value = 57756
$ws.Range("F194").Value = 2
$ws.Range("G194").Value = 2418
$ws.Range("B195").Value = 41737.86
$ws.Range("F225").Value = 17
$ws.Range("G225").Value = 1381.76
$ws.Range("B228").Value = 11659.59
$ws.Range("F313").Value = 72
$ws.Range("G313").Value = 15567.84
$ws.Range("F324").Value = 167
$ws.Range("G324").Value = 19076.41
$ws.Range("B380").Value = 255202.51
$ws.Range("F383").Value = 14
$ws.Range("G383").Value = 6253.1
$ws.Range("B389").Value = 22566.81
$ws.Range("B436").Value = 58047
$ws.Range("D436").Value = 105.54
$ws.Range("E436").Value = 126.1
$ws.Range("F436").Value = 62
$ws.Range("G436").Value = 6543.48
$ws.Range("B437").Value = 47097
$ws.Range("D437").Value = 112.28
$ws.Range("E437").Value = 134.16
$ws.Range("F437").Value = 15
$ws.Range("G437").Value = 1684.2
$ws.Range("F456").Value = 7
$ws.Range("G456").Value = 283.78
$ws.Range("F462").Value = 39
$ws.Range("G462").Value = 7271.16
$ws.Range("F466").Value = 46
$ws.Range("G466").Value = 2780.7
$ws.Range("B473").Value = 136376.12
$ws.Range("F491").Value = 473
$ws.Range("G491").Value = 6361.85
$ws.Range("F493").Value = 548
$ws.Range("G493").Value = 7019.88
$ws.Range("F494").Value = 262
$ws.Range("G494").Value = 6890.6
$ws.Range("F497").Value = 279
$ws.Range("G497").Value = 3573.99
$ws.Range("F501").Value = 88
$ws.Range("G501").Value = 1712.48
$ws.Range("F504").Value = 942
$ws.Range("G504").Value = 6113.58
$ws.Range("F506").Value = 322
$ws.Range("G506").Value = 8468.6
$ws.Range("F507").Value = 253
$ws.Range("G507").Value = 4156.79
$ws.Range("B509").Value = 94362.39
$ws.Range("F558").Value = 598
$ws.Range("G558").Value = 11870.3
$ws.Range("B563").Value = 36744.96
$ws.Range("F593").Value = 101
$ws.Range("G593").Value = 9948.5
$ws.Range("B601").Value = 64107.69
$ws.Range("F607").Value = 110
$ws.Range("G607").Value = 4987.4
$ws.Range("B613").Value = 6036.12
$ws.Range("F671").Value = 146
$ws.Range("G671").Value = 6304.28
$ws.Range("B677").Value = 20683.98
$ws.Range("F680").Value = 5
$ws.Range("G680").Value = 377.8
$ws.Range("F682").Value = 17
$ws.Range("G682").Value = 1394.68
$ws.Range("F683").Value = 25
$ws.Range("G683").Value = 2267
$ws.Range("F684").Value = 24
$ws.Range("G684").Value = 7461.12
$ws.Range("F687").Value = 33
$ws.Range("G687").Value = 3847.14
$ws.Range("F693").Value = 7
$ws.Range("G693").Value = 695.17
$ws.Range("B695").Value = 41803.95
$ws.Range("F755").Value = 227
$ws.Range("G755").Value = 18514.12
$ws.Range("F758").Value = 253
$ws.Range("G758").Value = 33016.5
$ws.Range("F763").Value = 106
$ws.Range("G763").Value = 2302.32
$ws.Range("F771").Value = 490
$ws.Range("G771").Value = 66154.89999999999
$ws.Range("F772").Value = 19
$ws.Range("G772").Value = 710.98
$ws.Range("F773").Value = 557
$ws.Range("G773").Value = 67235.47
$ws.Range("F774").Value = 43
$ws.Range("G774").Value = 5190.53
$ws.Range("B775").Value = 246677.63
$ws.Range("F784").Value = 181
$ws.Range("G784").Value = 5451.72
$ws.Range("B793").Value = 13373.21
$ws.Range("F800").Value = 7
$ws.Range("G800").Value = 261.8
$ws.Range("B801").Value = 455.66
$ws.Range("F852").Value = 624
$ws.Range("G852").Value = 18863.52
$ws.Range("F853").Value = 3152
$ws.Range("G853").Value = 514122.72
$ws.Range("F855").Value = 228
$ws.Range("G855").Value = 32980.2
$ws.Range("F856").Value = 120
$ws.Range("G856").Value = 4576.8
$ws.Range("F858").Value = 114
$ws.Range("G858").Value = 16856.04
$ws.Range("F860").Value = 127
$ws.Range("G860").Value = 16329.66
$ws.Range("B861").Value = 629625.98
$ws.Range("B867").Value = 3497502.56
$ws.Range("B868").Value = 3497502.56
